# TCD_PHIEU_TU_CHOI.docx edit script
# Implements:
#  1. "Ngày ${ngayThongBao}, ông (bà) "
#       -> "Ngày …. " + "tháng …. năm ……" + ", ông (bà) "   (3 runs, same rPr)
#  2. " đến Thanh tra thành phố để khiếu nại (tố cáo) về việc ${noiDung}."
#       -> " đến " + "${coQuanTiepNhan}" (bookmarked) + " để khiếu nại (tố cáo) về việc ${noiDung}."
#  3. "...,  Thanh tra thành phố nhận thấy..." -> "..., ${coQuanTiepNhan} nhận thấy..."
#  4. "...2013, Thanh tra thành phố từ chối..." -> "...2013, ${coQuanTiepNhan} từ chối..."
#  5. "Vậy, Thanh tra thành phố  thông báo..." -> "Vậy, ${coQuanTiepNhan} thông báo..."

$d = $word.ActiveDocument

function Find-Replace($findText, $replaceText) {
    $rng = $d.Content
    return $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

function Find-Start($findText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
    if (-not $ok) {
        return -1
    }
    return $rng.Start
}

function Split-At($pos) {
    # Insert then immediately delete a zero-width bookmark at a character
    # offset to force a run boundary there without altering the text or
    # merging runs back together afterwards.
    $r = $d.Range($pos, $pos)
    $name = "tmp_split_" + $pos
    $d.Bookmarks.Add($name, $r) | Out-Null
    $d.Bookmarks($name).Delete()
}

# ---------------------------------------------------------------------
# Change 1
# ---------------------------------------------------------------------
$old1 = 'Ngày ${ngayThongBao}, ông (bà) '
$seg1a = 'Ngày …. '
$seg1b = 'tháng …. năm ……'
$seg1c = ', ông (bà) '
$new1 = $seg1a + $seg1b + $seg1c

Find-Replace $old1 $new1 | Out-Null

$start1 = Find-Start $new1
if ($start1 -ge 0) {
    Split-At ($start1 + $seg1a.Length)
    Split-At ($start1 + $seg1a.Length + $seg1b.Length)
}

# ---------------------------------------------------------------------
# Change 2
# ---------------------------------------------------------------------
$old2 = ' đến Thanh tra thành phố để khiếu nại (tố cáo) về việc ${noiDung}.'
$seg2a = ' đến '
$seg2b = '${coQuanTiepNhan}'
$seg2c = ' để khiếu nại (tố cáo) về việc ${noiDung}.'
$new2 = $seg2a + $seg2b + $seg2c

Find-Replace $old2 $new2 | Out-Null

$start2 = Find-Start $new2
if ($start2 -ge 0) {
    $bmStart = $start2 + $seg2a.Length
    $bmEnd = $bmStart + $seg2b.Length
    $bmRange = $d.Range($bmStart, $bmEnd)
    $d.Bookmarks.Add('__DdeLink__1126_1750249842', $bmRange) | Out-Null
}

# ---------------------------------------------------------------------
# Change 3
# ---------------------------------------------------------------------
$old3 = 'Sau khi xem xét nội dung đơn khiếu nại (tố cáo) và nghe công dân trình bày, Thanh tra thành phố nhận thấy vụ việc đã được giải quyết đúng chính sách, pháp luật, được cơ quan Nhà nước có thẩm quyền kiểm tra, rà soát và thông báo trả lời bằng '
$new3 = 'Sau khi xem xét nội dung đơn khiếu nại (tố cáo) và nghe công dân trình bày, ${coQuanTiepNhan} nhận thấy vụ việc đã được giải quyết đúng chính sách, pháp luật, được cơ quan Nhà nước có thẩm quyền kiểm tra, rà soát và thông báo trả lời bằng '
Find-Replace $old3 $new3 | Out-Null

# ---------------------------------------------------------------------
# Change 4
# ---------------------------------------------------------------------
$old4 = 'Căn cứ Điều 9 Luật tiếp công dân năm 2013, Thanh tra thành phố từ chối tiếp nhận nội dung khiếu nại (tố cáo) của ông (bà) '
$new4 = 'Căn cứ Điều 9 Luật tiếp công dân năm 2013, ${coQuanTiepNhan} từ chối tiếp nhận nội dung khiếu nại (tố cáo) của ông (bà) '
Find-Replace $old4 $new4 | Out-Null

# ---------------------------------------------------------------------
# Change 5
# ---------------------------------------------------------------------
$old5 = 'Vậy, Thanh tra thành phố  thông báo để ông (bà) được biết; đề nghị ông (bà) chấp hành theo quy định của pháp luật, chấm dứt khiếu nại (tố cáo).'
$new5 = 'Vậy, ${coQuanTiepNhan} thông báo để ông (bà) được biết; đề nghị ông (bà) chấp hành theo quy định của pháp luật, chấm dứt khiếu nại (tố cáo).'
Find-Replace $old5 $new5 | Out-Null

Write-Host "done"
